$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression results with the new (peak-adjusted) values.
# Row 3
$ws.Range("C3").Value = [double]"1.03E-2"
$ws.Range("D3").Value = [double]"3.8300000000000001E-3"
$ws.Range("E3").Value = 2.68
$ws.Range("F3").Value = [double]"7.3000000000000001E-3"

# Row 4
$ws.Range("C4").Value = [double]"-2.1500000000000001E-5"
$ws.Range("D4").Value = [double]"9.91E-6"

# Row 5
$ws.Range("C5").Value = [double]"-7.3800000000000004E-2"
$ws.Range("C5").NumberFormat = "0.00E+00"

# Row 6
$ws.Range("C6").Value = 0.58799999999999997
$ws.Range("C6").NumberFormat = "0.00E+00"
$ws.Range("D6").Value = 0.217
$ws.Range("D6").NumberFormat = "0.00E+00"
$ws.Range("F6").Value = [double]"6.7000000000000002E-3"

# Row 7
$ws.Range("C7").Value = [double]"2.8299999999999999E-2"
$ws.Range("D7").Value = [double]"4.7499999999999999E-3"
$ws.Range("E7").Value = 5.95
$ws.Range("F7").Value = [double]"2.6000000000000001E-9"

# Row 8
$ws.Range("C8").Value = [double]"-4.1600000000000002E-5"
$ws.Range("D8").Value = [double]"1.01E-5"
$ws.Range("F8").Value = [double]"3.8000000000000002E-5"

# Row 9
$ws.Range("C9").Value = 0.61399999999999999
$ws.Range("C9").NumberFormat = "0.00E+00"
$ws.Range("D9").Value = 0.47399999999999998
$ws.Range("D9").NumberFormat = "0.00E+00"

# Row 10
$ws.Range("C10").Value = [double]"2.9700000000000001E-2"
$ws.Range("D10").Value = [double]"1.03E-2"
$ws.Range("E10").Value = 2.89
$ws.Range("F10").Value = [double]"3.8E-3"

# Row 11
$ws.Range("C11").Value = [double]"-1.4799999999999999E-4"
$ws.Range("D11").Value = [double]"5.7500000000000002E-5"
$ws.Range("F11").Value = [double]"1.03E-2"

# Row 13
$ws.Range("C13").Value = 0.02
$ws.Range("D13").Value = [double]"4.2700000000000004E-3"
$ws.Range("E13").Value = 4.6900000000000004
$ws.Range("F13").Value = [double]"2.7E-6"
$ws.Range("F13").NumberFormat = "0.00E+00"

# Row 14
$ws.Range("C14").Value = [double]"-4.6999999999999997E-5"
$ws.Range("D14").Value = [double]"1.19E-5"

# Update sheet view to match authored selection state.
[void]$ws.Range("E19").Select()
